# Updated simulation files with Holden scheme
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the now-unused duplicate trailing columns (X:AG), which held
#    copies of the N:W "pair" columns.
$ws.Range("X1:AG23").EntireColumn.Delete() | Out-Null

# 2. Re-order the HKL column headers in row 2 (columns C:M) to the new
#    sequence used by the updated simulation files.
$hklOrder = @("[2, 0, 0]", "[2, 2, 0]", "[3, 3, 3]", "[4, 2, 0]", "[4, 0, 0]", "[4, 2, 2]", "[5, 1, 1]", "[1, 1, 1]", "[2, 2, 2]", "[3, 3, 1]", "[3, 1, 1]")
for ($i = 0; $i -lt $hklOrder.Count; $i++) {
    $col = 3 + $i   # column C = 3
    $ws.Cells.Item(2, $col).Value = $hklOrder[$i]
}

# 3. Rows 16-19 (A=14..17) were labelled with the HexGrid scheme; they are
#    renamed to the new Holden scheme.
$holdenNames = @("Holden2.5", "Holden5", "Holden10", "Holden15")
for ($i = 0; $i -lt 4; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 2).Value = $holdenNames[$i]
}

# 4. Append new rows 20-23 carrying the HexGrid scheme (formerly on rows
#    16-19), with the same "all-ones" data pattern as the rest of the table.
$hexNames = @("HexGrid-90degTilt2.5degRes", "HexGrid-90degTilt5degRes", "HexGrid-90degTilt10degRes", "HexGrid-90degTilt15degRes")
for ($i = 0; $i -lt 4; $i++) {
    $row = 20 + $i
    $ws.Cells.Item($row, 1).Value = 18 + $i
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 4, 1).Style
    $ws.Cells.Item($row, 2).Value = $hexNames[$i]
    for ($col = 3; $col -le 23; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
